$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3766.8
$ws.Range("I64").Value = 3600.2222
$ws.Range("K64").Value = 3600.2222
$ws.Range("M64").Value = -3352.2222
$ws.Range("H67").Value = 3766.8
$ws.Range("I67").Value = 3600.2222
$ws.Range("K67").Value = 3600.2222
$ws.Range("M67").Value = -2742.2222
$ws.Range("H69").Value = 8562.375
$ws.Range("I69").Value = 1750
$ws.Range("K69").Value = 5250
$ws.Range("M69").Value = -4376
$ws.Range("H72").Value = 8562.375
$ws.Range("I72").Value = 1750
$ws.Range("K72").Value = 15750
$ws.Range("M72").Value = -11382
$ws.Range("H74").Value = 4485.143
$ws.Range("I74").Value = 4998
$ws.Range("K74").Value = 4998
$ws.Range("M74").Value = -4062
$ws.Range("H76").Value = 3037.3438
$ws.Range("I76").Value = 3027.4138
$ws.Range("K76").Value = 3027.4138
$ws.Range("M76").Value = -2712.4138
$ws.Range("H77").Value = 4485.143
$ws.Range("I77").Value = 4998
$ws.Range("K77").Value = 24990
$ws.Range("M77").Value = -20310
$ws.Range("H79").Value = 3037.3438
$ws.Range("I79").Value = 3027.4138
$ws.Range("K79").Value = 3027.4138
$ws.Range("M79").Value = -1935.4138
$ws.Range("H86").Value = 4832.8335
$ws.Range("I86").Value = 3416.2354
$ws.Range("K86").Value = 3416.2354
$ws.Range("M86").Value = -2293.2354
$ws.Range("H88").Value = 2935.7144
$ws.Range("I88").Value = 1998.4
$ws.Range("J88").Value = 3456.4443
$ws.Range("K88").Value = 1998.4
$ws.Range("L88").Value = 3456.4443
$ws.Range("M88").Value = -1592.4
$ws.Range("N88").Value = -4268.4443
$ws.Range("H89").Value = 4832.8335
$ws.Range("I89").Value = 3416.2354
$ws.Range("K89").Value = 17081.177
$ws.Range("M89").Value = -11465.177
$ws.Range("H91").Value = 2935.7144
$ws.Range("I91").Value = 1998.4
$ws.Range("J91").Value = 3456.4443
$ws.Range("K91").Value = 1998.4
$ws.Range("L91").Value = 3456.4443
$ws.Range("M91").Value = -594.4000000000001
$ws.Range("N91").Value = -6264.4443
$ws.Range("H112").Value = 1121.4286
$ws.Range("J112").Value = 1170.7693
$ws.Range("L112").Value = 3512.3079
$ws.Range("N112").Value = -5728.3079
$ws.Range("H137").Value = 2571.4849
$ws.Range("I137").Value = 1777.3478
$ws.Range("J137").Value = 4398
$ws.Range("K137").Value = 5332.0434
$ws.Range("L137").Value = 13194
$ws.Range("M137").Value = -2782.0434
$ws.Range("N137").Value = -18294

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3821.5293
$ws.Range("I32").Value = 3275.5112
$ws.Range("K32").Value = 3275.5112
$ws.Range("M32").Value = -2988.5112
$ws.Range("H45").Value = 1527.5
$ws.Range("I45").Value = 1036.2
$ws.Range("K45").Value = 1036.2
$ws.Range("M45").Value = -659.2
$ws.Range("H88").Value = 2385.9
$ws.Range("I88").Value = 2135.3333
$ws.Range("K88").Value = 2135.3333
$ws.Range("M88").Value = -1729.3333
$ws.Range("H91").Value = 2385.9
$ws.Range("I91").Value = 2135.3333
$ws.Range("K91").Value = 2135.3333
$ws.Range("M91").Value = -731.3332999999998
$ws.Range("H132").Value = 2343.353
$ws.Range("I132").Value = 2068.3076
$ws.Range("J132").Value = 3237.25
$ws.Range("K132").Value = 6204.9228
$ws.Range("L132").Value = 9711.75
$ws.Range("M132").Value = -3674.9228
$ws.Range("N132").Value = -14771.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1737.8125
$ws.Range("I86").Value = 1615.7693
$ws.Range("J86").Value = 2266.6667
$ws.Range("K86").Value = 1615.7693
$ws.Range("L86").Value = 2266.6667
$ws.Range("M86").Value = -492.7692999999999
$ws.Range("N86").Value = -4512.6667
$ws.Range("H89").Value = 1737.8125
$ws.Range("I89").Value = 1615.7693
$ws.Range("J89").Value = 2266.6667
$ws.Range("K89").Value = 8078.8465
$ws.Range("L89").Value = 11333.3335
$ws.Range("M89").Value = -2462.8465
$ws.Range("N89").Value = -22565.3335
$ws.Range("H99").Value = 1549.5714
$ws.Range("I99").Value = 1244.9445
$ws.Range("J99").Value = 2097.9
$ws.Range("K99").Value = 1244.9445
$ws.Range("L99").Value = 2097.9
$ws.Range("M99").Value = 253.0554999999999
$ws.Range("N99").Value = -5093.9
$ws.Range("H105").Value = 2302.7273
$ws.Range("I105").Value = 1880
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1880
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -133
$ws.Range("N105").Value = -5994
$ws.Range("H134").Value = 4214.2607
$ws.Range("I134").Value = 3960.4614
$ws.Range("J134").Value = 5628.2856
$ws.Range("K134").Value = 11881.3842
$ws.Range("L134").Value = 16884.8568
$ws.Range("M134").Value = -9346.3842
$ws.Range("N134").Value = -21954.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2856.625
$ws.Range("I62").Value = 2268.182
$ws.Range("J62").Value = 4151.2
$ws.Range("K62").Value = 2268.182
$ws.Range("L62").Value = 4151.2
$ws.Range("M62").Value = -1644.182
$ws.Range("N62").Value = -5399.2
$ws.Range("H65").Value = 2856.625
$ws.Range("I65").Value = 2268.182
$ws.Range("J65").Value = 4151.2
$ws.Range("K65").Value = 11340.91
$ws.Range("L65").Value = 20756
$ws.Range("M65").Value = -8220.91
$ws.Range("N65").Value = -26996
$ws.Range("H134").Value = 1504.4286
$ws.Range("I134").Value = 1457.6842
$ws.Range("J134").Value = 1948.5
$ws.Range("K134").Value = 4373.0526
$ws.Range("L134").Value = 5845.5
$ws.Range("M134").Value = -1838.0526
$ws.Range("N134").Value = -10915.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 520.4211
$ws.Range("I5").Value = 434.53333
$ws.Range("J5").Value = 842.5
$ws.Range("K5").Value = 1303.59999
$ws.Range("L5").Value = 2527.5
$ws.Range("M5").Value = -1191.59999
$ws.Range("N5").Value = -2751.5
$ws.Range("H132").Value = 3232
$ws.Range("I132").Value = 3093.3333
$ws.Range("J132").Value = 3786.6667
$ws.Range("K132").Value = 27839.9997
$ws.Range("L132").Value = 34080.0003
$ws.Range("M132").Value = -25309.9997
$ws.Range("N132").Value = -39140.0003
$ws.Range("H135").Value = 520.4211
$ws.Range("I135").Value = 434.53333
$ws.Range("J135").Value = 842.5
$ws.Range("K135").Value = 3910.79997
$ws.Range("L135").Value = 7582.5
$ws.Range("M135").Value = -1375.79997
$ws.Range("N135").Value = -12652.5
